$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.329.67"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "1.835.80"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.25%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5185"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3219"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06728"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7587"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.70%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07648"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.819.62"
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.000"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007855"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").Value = "26.337.20"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").Value = "2.076.16"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.535"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.379"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.905"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.220"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.651"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.156"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.122"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08704"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04765"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.857"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.111"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6881"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.051"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01752"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.184"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4808"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "110.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.70%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.077"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8804"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.589"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.91%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05845"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4104"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.960"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1226"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8783"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.76%  "
